# Update data: 5 February 2021
# Adds the newest (2021-01-01 / serial 44197) monthly observation to both
# the "Canada" sheet and the "Province" sheet.

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "Canada": append row 14
# -------------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$wsCanada.Cells.Item(14, 1).Value = 44197
$wsCanada.Cells.Item(14, 1).NumberFormat = "d-mmm-yy"

$wsCanada.Cells.Item(14, 2).Value = "Canada"
$wsCanada.Cells.Item(14, 2).NumberFormat = "d-mmm-yy"

$wsCanada.Cells.Item(14, 3).Value = 67
$wsCanada.Cells.Item(14, 4).Value = 1899

# Column A grows slightly wider once it best-fits the new date.
$wsCanada.Columns.Item(1).ColumnWidth = 9.498697916666666

# Move / record the selection the way the source workbook shows it.
$wsCanada.Range("C15").Select() | Out-Null

# -------------------------------------------------------------------------
# Sheet "Province": append rows 122-131 (one per province, same month)
# -------------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$provinceRows = @(
    @{ Row = 122; Name = "Newfoundland & Labrador"; C = 2.8;  D = 32.5 },
    @{ Row = 123; Name = "Prince Edward Island";     C = 0;    D = 6.7 },
    @{ Row = 124; Name = "Nova Scotia";              C = 10.5; D = 42.1 },
    @{ Row = 125; Name = "New Brunswick";            C = 14.5; D = 34.700000000000003 },
    @{ Row = 126; Name = "Quebec";                   C = 70.8; D = 394.6 },
    @{ Row = 127; Name = "Ontario";                  C = 93.1; D = 802.4 },
    @{ Row = 128; Name = "Manitoba";                 C = 51.4; D = 54.8 },
    @{ Row = 129; Name = "Saskatchewan";              C = 13.7; D = 42.3 },
    @{ Row = 130; Name = "Alberta";                  C = 44.4; D = 262.7 },
    @{ Row = 131; Name = "British Columbia";         C = 76.2; D = 226.2 }
)

foreach ($r in $provinceRows) {
    $row = $r.Row
    $wsProvince.Cells.Item($row, 1).Value = 44197
    $wsProvince.Cells.Item($row, 1).NumberFormat = "d-mmm-yy"

    $wsProvince.Cells.Item($row, 2).Value = $r.Name

    $wsProvince.Cells.Item($row, 3).Value = $r.C
    $wsProvince.Cells.Item($row, 4).Value = $r.D
}

# First data row of the new block keeps the date's style (like row 122 in the
# diff, which also carries the style on column B).
$wsProvince.Cells.Item(122, 2).NumberFormat = "d-mmm-yy"

$wsProvince.Activate() | Out-Null
$wsProvince.Range("C132").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 116
$excel.ActiveWindow.ScrollColumn = 1
